$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# block1
# row 5 <- row 10
$ws.Range("A5").Value2 = 111621552
$ws.Range("B5").Value2 = 56398
$ws.Range("D5").Value2 = "NT"
$ws.Range("E5").Value2 = 100109
$ws.Range("F5").Value2 = "Tretåig hackspett"
$ws.Range("G5").Value2 = "Picoides tridactylus"
$ws.Range("H5").Value2 = "(Linnaeus, 1758)"
$ws.Range("M5").Value2 = "färska spår"
$ws.Range("Q5").Value2 = 536390.751010091
$ws.Range("R5").Value2 = 7208948.130241925
$ws.Range("S5").Value2 = 10
$ws.Range("AC5").Value2 = $null
$ws.Range("AJ5").Value2 = "gran"
$ws.Range("AK5").Value2 = "Picea abies"
$ws.Range("AO5").Value2 = "Picea abies"

# row 6 <- row 5
$ws.Range("A6").Value2 = 111621565
$ws.Range("B6").Value2 = 56398
$ws.Range("D6").Value2 = "NT"
$ws.Range("E6").Value2 = 100109
$ws.Range("F6").Value2 = "Tretåig hackspett"
$ws.Range("G6").Value2 = "Picoides tridactylus"
$ws.Range("H6").Value2 = "(Linnaeus, 1758)"
$ws.Range("M6").Value2 = "äldre spår"
$ws.Range("Q6").Value2 = 536384.516595734
$ws.Range("R6").Value2 = 7208973.017290085
$ws.Range("S6").Value2 = 25
$ws.Range("AC6").Value2 = $null
$ws.Range("AJ6").Value2 = "gran"
$ws.Range("AK6").Value2 = "Picea abies"
$ws.Range("AO6").Value2 = "Picea abies"

# row 7 <- row 6
$ws.Range("A7").Value2 = 111621572
$ws.Range("B7").Value2 = 90087
$ws.Range("D7").Value2 = "LC"
$ws.Range("E7").Value2 = 3298
$ws.Range("F7").Value2 = "Trådticka"
$ws.Range("G7").Value2 = "Climacocystis borealis"
$ws.Range("H7").Value2 = "(Fr.) Kotl. & Pouzar"
$ws.Range("M7").Value2 = $null
$ws.Range("Q7").Value2 = 536384.516595734
$ws.Range("R7").Value2 = 7208973.017290085
$ws.Range("S7").Value2 = 25
$ws.Range("AC7").Value2 = "Vid stambasen, stående gran"
$ws.Range("AJ7").Value2 = "gran"
$ws.Range("AK7").Value2 = "Picea abies"
$ws.Range("AO7").Value2 = "Picea abies"

# row 8 <- row 7
$ws.Range("A8").Value2 = 111624097
$ws.Range("B8").Value2 = 56398
$ws.Range("D8").Value2 = "NT"
$ws.Range("E8").Value2 = 100109
$ws.Range("F8").Value2 = "Tretåig hackspett"
$ws.Range("G8").Value2 = "Picoides tridactylus"
$ws.Range("H8").Value2 = "(Linnaeus, 1758)"
$ws.Range("M8").Value2 = "äldre spår"
$ws.Range("Q8").Value2 = 535928.2260010642
$ws.Range("R8").Value2 = 7209575.417457776
$ws.Range("S8").Value2 = 25
$ws.Range("AC8").Value2 = $null
$ws.Range("AJ8").Value2 = $null
$ws.Range("AK8").Value2 = $null
$ws.Range("AO8").Value2 = $null

# row 9 <- row 8
$ws.Range("A9").Value2 = 111623672
$ws.Range("B9").Value2 = 90087
$ws.Range("D9").Value2 = "LC"
$ws.Range("E9").Value2 = 3298
$ws.Range("F9").Value2 = "Trådticka"
$ws.Range("G9").Value2 = "Climacocystis borealis"
$ws.Range("H9").Value2 = "(Fr.) Kotl. & Pouzar"
$ws.Range("M9").Value2 = $null
$ws.Range("Q9").Value2 = 535988.4578048707
$ws.Range("R9").Value2 = 7209708.122271948
$ws.Range("S9").Value2 = 25
$ws.Range("AC9").Value2 = $null
$ws.Range("AJ9").Value2 = $null
$ws.Range("AK9").Value2 = $null
$ws.Range("AO9").Value2 = $null

# row 10 <- row 9
$ws.Range("A10").Value2 = 111625174
$ws.Range("B10").Value2 = 56398
$ws.Range("D10").Value2 = "NT"
$ws.Range("E10").Value2 = 100109
$ws.Range("F10").Value2 = "Tretåig hackspett"
$ws.Range("G10").Value2 = "Picoides tridactylus"
$ws.Range("H10").Value2 = "(Linnaeus, 1758)"
$ws.Range("M10").Value2 = "äldre spår"
$ws.Range("Q10").Value2 = 536263.4947354996
$ws.Range("R10").Value2 = 7209449.609840255
$ws.Range("S10").Value2 = 25
$ws.Range("AC10").Value2 = "Ringhack på 2 granar"
$ws.Range("AJ10").Value2 = $null
$ws.Range("AK10").Value2 = $null
$ws.Range("AO10").Value2 = $null

# block2
# row 50 <- row 51
$ws.Range("A50").Value2 = 111622736
$ws.Range("B50").Value2 = 90087
$ws.Range("D50").Value2 = "LC"
$ws.Range("E50").Value2 = 3298
$ws.Range("F50").Value2 = "Trådticka"
$ws.Range("G50").Value2 = "Climacocystis borealis"
$ws.Range("H50").Value2 = "(Fr.) Kotl. & Pouzar"
$ws.Range("M50").Value2 = $null
$ws.Range("Q50").Value2 = 535953.8130829642
$ws.Range("R50").Value2 = 7209209.795134133
$ws.Range("S50").Value2 = 10

# row 51 <- row 52
$ws.Range("A51").Value2 = 111623417
$ws.Range("B51").Value2 = 90087
$ws.Range("D51").Value2 = "LC"
$ws.Range("E51").Value2 = 3298
$ws.Range("F51").Value2 = "Trådticka"
$ws.Range("G51").Value2 = "Climacocystis borealis"
$ws.Range("H51").Value2 = "(Fr.) Kotl. & Pouzar"
$ws.Range("M51").Value2 = $null
$ws.Range("Q51").Value2 = 535866.1958485778
$ws.Range("R51").Value2 = 7209556.480484258
$ws.Range("S51").Value2 = 25

# row 52 <- row 53
$ws.Range("A52").Value2 = 111624796
$ws.Range("B52").Value2 = 89686
$ws.Range("D52").Value2 = "NT"
$ws.Range("E52").Value2 = 658
$ws.Range("F52").Value2 = "Rosenticka"
$ws.Range("G52").Value2 = "Rhodofomes roseus"
$ws.Range("H52").Value2 = "(Alb. & Schwein.) Kotl. & Pouzar"
$ws.Range("M52").Value2 = $null
$ws.Range("Q52").Value2 = 536163.445390123
$ws.Range("R52").Value2 = 7209387.476444452
$ws.Range("S52").Value2 = 25

# row 53 <- row 54
$ws.Range("A53").Value2 = 111623380
$ws.Range("B53").Value2 = 56398
$ws.Range("D53").Value2 = "NT"
$ws.Range("E53").Value2 = 100109
$ws.Range("F53").Value2 = "Tretåig hackspett"
$ws.Range("G53").Value2 = "Picoides tridactylus"
$ws.Range("H53").Value2 = "(Linnaeus, 1758)"
$ws.Range("M53").Value2 = "äldre spår"
$ws.Range("Q53").Value2 = 535912.0367731415
$ws.Range("R53").Value2 = 7209513.038373807
$ws.Range("S53").Value2 = 25

# row 54 <- row 50
$ws.Range("A54").Value2 = 111622993
$ws.Range("B54").Value2 = 56398
$ws.Range("D54").Value2 = "NT"
$ws.Range("E54").Value2 = 100109
$ws.Range("F54").Value2 = "Tretåig hackspett"
$ws.Range("G54").Value2 = "Picoides tridactylus"
$ws.Range("H54").Value2 = "(Linnaeus, 1758)"
$ws.Range("M54").Value2 = "äldre spår"
$ws.Range("Q54").Value2 = 536109.6063802312
$ws.Range("R54").Value2 = 7209286.560724956
$ws.Range("S54").Value2 = 50

# block3
# row 62 <- row 66
$ws.Range("A62").Value2 = 111624781
$ws.Range("B62").Value2 = 89405
$ws.Range("E62").Value2 = 1202
$ws.Range("F62").Value2 = "Ullticka"
$ws.Range("G62").Value2 = "Phellinidium ferrugineofuscum"
$ws.Range("H62").Value2 = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("M62").Value2 = $null
$ws.Range("Q62").Value2 = 536163.445390123
$ws.Range("R62").Value2 = 7209387.476444452
$ws.Range("S62").Value2 = 25
$ws.Range("AC62").Value2 = "Även rosenticka mfl vedsvampar på denna låga"

# row 63 <- row 62
$ws.Range("A63").Value2 = 111622734
$ws.Range("B63").Value2 = 56398
$ws.Range("E63").Value2 = 100109
$ws.Range("F63").Value2 = "Tretåig hackspett"
$ws.Range("G63").Value2 = "Picoides tridactylus"
$ws.Range("H63").Value2 = "(Linnaeus, 1758)"
$ws.Range("M63").Value2 = "äldre spår"
$ws.Range("Q63").Value2 = 535953.8130829642
$ws.Range("R63").Value2 = 7209209.795134133
$ws.Range("S63").Value2 = 10
$ws.Range("AC63").Value2 = $null

# row 64 <- row 63
$ws.Range("A64").Value2 = 111623737
$ws.Range("B64").Value2 = 56398
$ws.Range("E64").Value2 = 100109
$ws.Range("F64").Value2 = "Tretåig hackspett"
$ws.Range("G64").Value2 = "Picoides tridactylus"
$ws.Range("H64").Value2 = "(Linnaeus, 1758)"
$ws.Range("M64").Value2 = "färska spår"
$ws.Range("Q64").Value2 = 535968.9484369244
$ws.Range("R64").Value2 = 7209745.533198988
$ws.Range("S64").Value2 = 25
$ws.Range("AC64").Value2 = $null

# row 65 <- row 64
$ws.Range("A65").Value2 = 111624558
$ws.Range("B65").Value2 = 89790
$ws.Range("E65").Value2 = 6040186
$ws.Range("F65").Value2 = $null
$ws.Range("G65").Value2 = "Leptoporus mollis"
$ws.Range("H65").Value2 = "(Pers.:Fr.) Quél."
$ws.Range("M65").Value2 = $null
$ws.Range("Q65").Value2 = 536083.1087774199
$ws.Range("R65").Value2 = 7209411.039029445
$ws.Range("S65").Value2 = 25
$ws.Range("AC65").Value2 = $null

# row 66 <- row 65
$ws.Range("A66").Value2 = 111625227
$ws.Range("B66").Value2 = 77515
$ws.Range("E66").Value2 = 6425
$ws.Range("F66").Value2 = "Garnlav"
$ws.Range("G66").Value2 = "Alectoria sarmentosa"
$ws.Range("H66").Value2 = "(Ach.) Ach."
$ws.Range("M66").Value2 = $null
$ws.Range("Q66").Value2 = 536368.8900330348
$ws.Range("R66").Value2 = 7209489.813207326
$ws.Range("S66").Value2 = 25
$ws.Range("AC66").Value2 = $null
